$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 12000
$ws.Range("M46").Value = -11881

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 4000
$ws.Range("I60").Value = 4000
$ws.Range("K60").Value = 12000
$ws.Range("M60").Value = -11516

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 4036.3333
$ws.Range("I111").Value = 5304.5
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 15913.5
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = -12846.5
$ws.Range("N111").Value = -10634

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4829.0933
$ws.Range("J112").Value = 5285.795
$ws.Range("L112").Value = 15857.385
$ws.Range("N112").Value = -18073.385

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3363.182
$ws.Range("I116").Value = 3374.375
$ws.Range("J116").Value = 3333.3333
$ws.Range("K116").Value = 3374.375
$ws.Range("L116").Value = 3333.3333
$ws.Range("M116").Value = 67.625
$ws.Range("N116").Value = -10217.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 983.8431399999999
$ws.Range("I129").Value = 427.27274
$ws.Range("J129").Value = 1136.9
$ws.Range("K129").Value = 1281.81822
$ws.Range("L129").Value = 3410.7
$ws.Range("M129").Value = 3718.18178
$ws.Range("N129").Value = -13410.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 13514454
$ws.Range("I135").Value = 16129573
$ws.Range("K135").Value = 145166157
$ws.Range("M135").Value = -145163622

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2877.83
$ws.Range("I138").Value = 1504.625
$ws.Range("J138").Value = 3139.3928
$ws.Range("K138").Value = 4513.875
$ws.Range("L138").Value = 9418.178400000001
$ws.Range("M138").Value = 626.125
$ws.Range("N138").Value = -19698.1784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5557596
$ws.Range("I61").Value = 7753384.5
$ws.Range("J61").Value = 3543.353
$ws.Range("K61").Value = 7753384.5
$ws.Range("L61").Value = 3543.353
$ws.Range("M61").Value = -7753172.5
$ws.Range("N61").Value = -3967.353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 43618.2
$ws.Range("J64").Value = 43618.2
$ws.Range("L64").Value = 43618.2
$ws.Range("N64").Value = -44114.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 43618.2
$ws.Range("J67").Value = 43618.2
$ws.Range("L67").Value = 43618.2
$ws.Range("N67").Value = -45334.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2340.8293
$ws.Range("I74").Value = 1527.2142
$ws.Range("J74").Value = 4093.2307
$ws.Range("K74").Value = 1527.2142
$ws.Range("L74").Value = 4093.2307
$ws.Range("M74").Value = -653.2141999999999
$ws.Range("N74").Value = -5841.2307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2340.8293
$ws.Range("I77").Value = 1527.2142
$ws.Range("J77").Value = 4093.2307
$ws.Range("K77").Value = 7636.071
$ws.Range("L77").Value = 20466.1535
$ws.Range("M77").Value = -3268.071
$ws.Range("N77").Value = -29202.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 66000
$ws.Range("J92").Value = 66000
$ws.Range("L92").Value = 66000
$ws.Range("N92").Value = -70992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4258.7
$ws.Range("I132").Value = 3807.9644
$ws.Range("J132").Value = 5310.4165
$ws.Range("K132").Value = 11423.8932
$ws.Range("L132").Value = 15931.2495
$ws.Range("M132").Value = -8893.893199999999
$ws.Range("N132").Value = -20991.2495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5557596
$ws.Range("I136").Value = 7753384.5
$ws.Range("J136").Value = 3543.353
$ws.Range("K136").Value = 23260153.5
$ws.Range("L136").Value = 10630.059
$ws.Range("M136").Value = -23257603.5
$ws.Range("N136").Value = -15730.059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 53024
$ws.Range("J132").Value = 53780
$ws.Range("L132").Value = 53780
$ws.Range("N132").Value = -63900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3063.0435
$ws.Range("I134").Value = 3217.8333
$ws.Range("J134").Value = 2894.182
$ws.Range("K134").Value = 9653.499899999999
$ws.Range("L134").Value = 8682.545999999998
$ws.Range("M134").Value = -7118.499899999999
$ws.Range("N134").Value = -13752.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5866.45
$ws.Range("I31").Value = 1891.0555
$ws.Range("J31").Value = 7570.1904
$ws.Range("K31").Value = 1891.0555
$ws.Range("L31").Value = 7570.1904
$ws.Range("M31").Value = -1596.0555
$ws.Range("N31").Value = -8160.1904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5866.45
$ws.Range("I34").Value = 1891.0555
$ws.Range("J34").Value = 7570.1904
$ws.Range("K34").Value = 1891.0555
$ws.Range("L34").Value = 7570.1904
$ws.Range("M34").Value = -1689.0555
$ws.Range("N34").Value = -7974.1904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 13316.077
$ws.Range("I60").Value = 850
$ws.Range("J60").Value = 15582.637
$ws.Range("K60").Value = 850
$ws.Range("L60").Value = 15582.637
$ws.Range("M60").Value = -339
$ws.Range("N60").Value = -16604.637

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1428
$ws.Range("I68").Value = 1250.5
$ws.Range("J68").Value = 1605.5
$ws.Range("K68").Value = 3751.5
$ws.Range("L68").Value = 4816.5
$ws.Range("M68").Value = -2940.5
$ws.Range("N68").Value = -6438.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1428
$ws.Range("I71").Value = 1250.5
$ws.Range("J71").Value = 1605.5
$ws.Range("K71").Value = 11254.5
$ws.Range("L71").Value = 14449.5
$ws.Range("M71").Value = -7198.5
$ws.Range("N71").Value = -22561.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 879.8
$ws.Range("I86").Value = 799.6667
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 2399.0001
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1213.0001
$ws.Range("N86").Value = -5372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 879.8
$ws.Range("I89").Value = 799.6667
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 7197.0003
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -1269.0003
$ws.Range("N89").Value = -20856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 4000
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -16118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2626
$ws.Range("J107").Value = 3660.4
$ws.Range("L107").Value = 10981.2
$ws.Range("N107").Value = -14821.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 846.75
$ws.Range("J113").Value = 811.3333
$ws.Range("L113").Value = 2433.9999
$ws.Range("N113").Value = -6773.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3063.4666
$ws.Range("I122").Value = 495.82352
$ws.Range("J122").Value = 4622.393
$ws.Range("K122").Value = 4462.41168
$ws.Range("L122").Value = 41601.537
$ws.Range("M122").Value = -2012.41168
$ws.Range("N122").Value = -46501.537

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4052.6667
$ws.Range("I131").Value = 766
$ws.Range("J131").Value = 4536
$ws.Range("K131").Value = 2298
$ws.Range("L131").Value = 13608
$ws.Range("M131").Value = 2742
$ws.Range("N131").Value = -23688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2302.6
$ws.Range("I136").Value = 853.25
$ws.Range("J136").Value = 8100
$ws.Range("K136").Value = 2559.75
$ws.Range("L136").Value = 24300
$ws.Range("M136").Value = 2540.25
$ws.Range("N136").Value = -34500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 47726.832
$ws.Range("J137").Value = 146271.42
$ws.Range("L137").Value = 438814.26
$ws.Range("N137").Value = -449014.26

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1763.1471
$ws.Range("I140").Value = 1160.35
$ws.Range("J140").Value = 2624.2856
$ws.Range("K140").Value = 3481.05
$ws.Range("L140").Value = 7872.8568
$ws.Range("M140").Value = 1698.95
$ws.Range("N140").Value = -18232.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 7667.222
$ws.Range("I141").Value = 3876.25
$ws.Range("J141").Value = 10700
$ws.Range("K141").Value = 11628.75
$ws.Range("L141").Value = 32100
$ws.Range("M141").Value = -6448.75
$ws.Range("N141").Value = -42460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 14063.611
$ws.Range("J123").Value = 17224.643
$ws.Range("L123").Value = 17224.643
$ws.Range("N123").Value = -22124.643

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 38464430
$ws.Range("I132").Value = 58825716
$ws.Range("J132").Value = 4224.8887
$ws.Range("K132").Value = 176477148
$ws.Range("L132").Value = 12674.6661
$ws.Range("M132").Value = -176474618
$ws.Range("N132").Value = -17734.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2950
$ws.Range("I46").Value = 675
$ws.Range("K46").Value = 675
$ws.Range("M46").Value = -487

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1921.8334
$ws.Range("I82").Value = 1400.2858
$ws.Range("J82").Value = 2253.7273
$ws.Range("K82").Value = 1400.2858
$ws.Range("L82").Value = 2253.7273
$ws.Range("M82").Value = -1039.2858
$ws.Range("N82").Value = -2975.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1921.8334
$ws.Range("I85").Value = 1400.2858
$ws.Range("J85").Value = 2253.7273
$ws.Range("K85").Value = 1400.2858
$ws.Range("L85").Value = 2253.7273
$ws.Range("M85").Value = -152.2858000000001
$ws.Range("N85").Value = -4749.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1948
$ws.Range("I136").Value = 1705
$ws.Range("J136").Value = 2579.8
$ws.Range("K136").Value = 5115
$ws.Range("L136").Value = 7739.400000000001
$ws.Range("M136").Value = -2565
$ws.Range("N136").Value = -12839.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 55726.082
$ws.Range("J140").Value = 55726.082
$ws.Range("L140").Value = 55726.082
$ws.Range("N140").Value = -66086.08199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 34990
$ws.Range("J99").Value = 34990
$ws.Range("L99").Value = 34990
$ws.Range("N99").Value = -40980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2545.3
$ws.Range("I122").Value = 2633.6667
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 7901.000100000001
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -5451.000100000001
$ws.Range("N122").Value = -10150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4905223.5
$ws.Range("I132").Value = 3269.4546
$ws.Range("K132").Value = 9808.363799999999
$ws.Range("M132").Value = -7278.363799999999
